$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 294, shifting existing rows 294:428 down to 295:429
$ws.Rows.Item(294).Insert()

# Populate the newly inserted row 294 with the new record
$ws.Cells.Item(294, 1).Value2  = 3
$ws.Cells.Item(294, 2).Value2  = "Femacal de La Calera"
$ws.Cells.Item(294, 3).Value2  = "Coquimbo"
$ws.Cells.Item(294, 4).Value2  = 44609
$ws.Cells.Item(294, 5).Value2  = 5
$ws.Cells.Item(294, 6).Value2  = 100112037
$ws.Cells.Item(294, 7).Value2  = "Cebollín"
$ws.Cells.Item(294, 8).Value2  = "Sin especificar"
$ws.Cells.Item(294, 9).Value2  = "Primera"
$ws.Cells.Item(294, 10).Value2 = 180
$ws.Cells.Item(294, 11).Value2 = 3500
$ws.Cells.Item(294, 12).Value2 = 3500
$ws.Cells.Item(294, 13).Value2 = 3500
$ws.Cells.Item(294, 14).Value2 = "$/paquete 36 unidades"
$ws.Cells.Item(294, 15).Value2 = "Provincia de Quillota"
$ws.Cells.Item(294, 16).Value2 = 97
$ws.Cells.Item(294, 17).Value2 = 36
$ws.Cells.Item(294, 18).Value2 = "Hortaliza"

# Ensure the date cell keeps the same style as the rest of column D (style index 2 / datetime format)
$ws.Cells.Item(294, 4).NumberFormat = $ws.Cells.Item(295, 4).NumberFormat
